$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 500
$ws.Range("I12").Value = 500
$ws.Range("K12").Value = 500
$ws.Range("M12").Value = -330

$ws.Range("H33").Value = 241.10527
$ws.Range("I33").Value = 113.882355
$ws.Range("K33").Value = 113.882355
$ws.Range("M33").Value = 115.117645

$ws.Range("H39").Value = 203
$ws.Range("I39").Value = 71.666664
$ws.Range("K39").Value = 214.999992
$ws.Range("M39").Value = 81.00000800000001

$ws.Range("H62").Value = 5603.143
$ws.Range("I62").Value = 3257.3333
$ws.Range("J62").Value = 7362.5
$ws.Range("K62").Value = 3257.3333
$ws.Range("L62").Value = 7362.5
$ws.Range("M62").Value = -2633.3333
$ws.Range("N62").Value = -8610.5

$ws.Range("H65").Value = 5603.143
$ws.Range("I65").Value = 3257.3333
$ws.Range("J65").Value = 7362.5
$ws.Range("K65").Value = 16286.6665
$ws.Range("L65").Value = 36812.5
$ws.Range("M65").Value = -13166.6665
$ws.Range("N65").Value = -43052.5

$ws.Range("H92").Value = 983.2143
$ws.Range("I92").Value = 987.2727
$ws.Range("J92").Value = 968.3333
$ws.Range("K92").Value = 987.2727
$ws.Range("L92").Value = 968.3333
$ws.Range("M92").Value = 260.7273
$ws.Range("N92").Value = -3464.3333

$ws.Range("H116").Value = 6193.75
$ws.Range("I116").Value = 3390
$ws.Range("K116").Value = 3390
$ws.Range("M116").Value = 52

$ws.Range("H138").Value = 6783.077
$ws.Range("J138").Value = 5407.4165
$ws.Range("L138").Value = 16222.2495
$ws.Range("N138").Value = -26502.2495

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 15223.139
$ws.Range("I32").Value = 6342.952
$ws.Range("K32").Value = 6342.952
$ws.Range("M32").Value = -6055.952

$ws.Range("H45").Value = 2422.111
$ws.Range("I45").Value = 800
$ws.Range("K45").Value = 800
$ws.Range("M45").Value = -423

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 720.4286
$ws.Range("I94").Value = 648.6
$ws.Range("K94").Value = 648.6
$ws.Range("M94").Value = -197.6

$ws.Range("H134").Value = 2062.6572
$ws.Range("I134").Value = 1603.4517
$ws.Range("J134").Value = 5621.5
$ws.Range("K134").Value = 4810.355100000001
$ws.Range("L134").Value = 16864.5
$ws.Range("M134").Value = -2275.355100000001
$ws.Range("N134").Value = -21934.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 266.33334
$ws.Range("I22").Value = 149.33333
$ws.Range("J22").Value = 383.33334
$ws.Range("K22").Value = 149.33333
$ws.Range("L22").Value = 383.33334
$ws.Range("M22").Value = 200.66667
$ws.Range("N22").Value = -1083.33334

$ws.Range("H62").Value = 48869.777
$ws.Range("I62").Value = 4964.75
$ws.Range("K62").Value = 4964.75
$ws.Range("M62").Value = -4340.75

$ws.Range("H65").Value = 48869.777
$ws.Range("I65").Value = 4964.75
$ws.Range("K65").Value = 24823.75
$ws.Range("M65").Value = -21703.75

$ws.Range("H132").Value = 856.4706
$ws.Range("I132").Value = 856.4706
$ws.Range("K132").Value = 2569.4118
$ws.Range("M132").Value = -39.41179999999986

$ws.Range("H134").Value = 2680.3635
$ws.Range("I134").Value = 1971.5
$ws.Range("K134").Value = 5914.5
$ws.Range("M134").Value = -3379.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 2009.5555
$ws.Range("I131").Value = 772.75
$ws.Range("K131").Value = 2318.25
$ws.Range("M131").Value = 2721.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 10520000
$ws.Range("J18").Value = 40000
$ws.Range("L18").Value = 40000
$ws.Range("N18").Value = -40586

$ws.Range("H29").Value = 18916.666
$ws.Range("I29").Value = 19000
$ws.Range("J29").Value = 18888.889
$ws.Range("K29").Value = 19000
$ws.Range("L29").Value = 18888.889
$ws.Range("M29").Value = -18710
$ws.Range("N29").Value = -19468.889

$ws.Range("H80").Value = 8441.25
$ws.Range("J80").Value = 8921.666999999999
$ws.Range("L80").Value = 8921.666999999999
$ws.Range("N80").Value = -10917.667

$ws.Range("H83").Value = 8441.25
$ws.Range("J83").Value = 8921.666999999999
$ws.Range("L83").Value = 44608.335
$ws.Range("N83").Value = -54592.335

$ws.Range("H97").Value = 1399.8148
$ws.Range("I97").Value = 1521
$ws.Range("K97").Value = 1521
$ws.Range("M97").Value = -1025

$ws.Range("H102").Value = 2446.4443
$ws.Range("I102").Value = 1485.091
$ws.Range("K102").Value = 1485.091
$ws.Range("M102").Value = 136.9090000000001

$ws.Range("H132").Value = 2837.0833
$ws.Range("I132").Value = 2287
$ws.Range("K132").Value = 6861
$ws.Range("M132").Value = -4331

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1176.6
$ws.Range("I22").Value = 703.5
$ws.Range("J22").Value = 1492
$ws.Range("K22").Value = 703.5
$ws.Range("L22").Value = 1492
$ws.Range("M22").Value = -408.5
$ws.Range("N22").Value = -2082

$ws.Range("H27").Value = 1176.6
$ws.Range("I27").Value = 703.5
$ws.Range("J27").Value = 1492
$ws.Range("K27").Value = 703.5
$ws.Range("L27").Value = 1492
$ws.Range("M27").Value = -596.5
$ws.Range("N27").Value = -1706

$ws.Range("H82").Value = 2826.111
$ws.Range("I82").Value = 2554.375
$ws.Range("J82").Value = 5000
$ws.Range("K82").Value = 2554.375
$ws.Range("L82").Value = 5000
$ws.Range("M82").Value = -2193.375
$ws.Range("N82").Value = -5722

$ws.Range("H85").Value = 2826.111
$ws.Range("I85").Value = 2554.375
$ws.Range("J85").Value = 5000
$ws.Range("K85").Value = 2554.375
$ws.Range("L85").Value = 5000
$ws.Range("M85").Value = -1306.375
$ws.Range("N85").Value = -7496

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H19").Value = 5001500
$ws.Range("J19").Value = 3000
$ws.Range("L19").Value = 3000
$ws.Range("N19").Value = -3348

$ws.Range("H46").Value = 78795.836
$ws.Range("J46").Value = 78795.836
$ws.Range("L46").Value = 78795.836
$ws.Range("N46").Value = -79257.836

$ws.Range("H113").Value = 1404.9333
$ws.Range("I113").Value = 834.36365
$ws.Range("J113").Value = 2974
$ws.Range("K113").Value = 2503.09095
$ws.Range("L113").Value = 8922
$ws.Range("M113").Value = -333.0909499999998
$ws.Range("N113").Value = -13262

$ws.Range("H134").Value = 78795.836
$ws.Range("J134").Value = 78795.836
$ws.Range("L134").Value = 236387.508
$ws.Range("N134").Value = -241457.508

$ws.Range("H136").Value = 44085.582
$ws.Range("I136").Value = 1392.1177
$ws.Range("K136").Value = 4176.3531
$ws.Range("M136").Value = -1626.3531
